$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) Pewter City paragraph: fix "Terra" (grammar-checker) artifact -
#    collapse ", Quick Balls somewhere, " + "Cannot" + " leave the city
#    until you defeat Brock." into a single run of plain text, which
#    also drops the now-stale proofErr gramStart/gramEnd markers.
# ---------------------------------------------------------------------
$fixFind = $d.Content.Find
$fixFind.ClearFormatting()
$fixFind.Replacement.ClearFormatting()
$fixFind.Text = ", Quick Balls somewhere, Cannot leave the city until you defeat Brock."
$fixFind.Replacement.Text = ", Quick Balls somewhere, Cannot leave the city until you defeat Brock."
$fixFind.Execute($null, $false, $false, $false, $false, $false, $true, 1, $false, $null, 2) | Out-Null

# ---------------------------------------------------------------------
# 2) Pewter Gym trainers added: insert a bold " (ADDED) " right after
#    the bold "Pewter City" run.
# ---------------------------------------------------------------------
$target = $d.Content
$target.Find.ClearFormatting()
$target.Find.Execute("Pewter City", $false, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$target.Collapse(0)
$target.InsertAfter(" (ADDED) ")
$target.Font.Bold = 1

# ---------------------------------------------------------------------
# 3) Move the "_GoBack" bookmark from the end of the document (after
#    "Saffron City-") up to right after the newly-added " (ADDED) "
#    text. Word only ever keeps one "_GoBack" bookmark, so re-adding it
#    here automatically removes the old one further down.
# ---------------------------------------------------------------------
$target.Collapse(0)
$d.Bookmarks.Add("_GoBack", $target)
